# Lesson 10 slide deck update:
#  1. Merge split text runs on slide 1 ("Demo "/"Lab 1..." and "Lab "/"2 prelab...")
#  2. Merge split text runs on slide 11 ("Jeff Falkinburg" / ", ")
#  3. Insert a new "Application Binary Interface (ABI)" slide just before the
#     final "Lab 2 Introduction" slide (i.e. as the new second-to-last slide).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 1 - "Lesson Outline" bullet list: collapse two runs that were
#    split mid-sentence back into a single run each.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$outline = $slide1.Shapes.Item(2).TextFrame.TextRange

for ($i = 1; $i -le $outline.Paragraphs().Count; $i++) {
    $para = $outline.Paragraphs($i)
    $t = $para.Text.TrimEnd([char]13, [char]10)
    if ($t -eq "Demo Lab 1 functionality by COB today!") {
        $para.Text = "temp_merge_placeholder"
        $para.Text = "Demo Lab 1 functionality by COB today!"
    } elseif ($t -eq "Lab 2 prelab due BOC next lesson") {
        $para.Text = "temp_merge_placeholder"
        $para.Text = "Lab 2 prelab due BOC next lesson"
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 11 - code comment line: merge "Jeff Falkinburg" + ", " into one
#    run while leaving the trailing "USAF" run untouched.
# ---------------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$code = $slide11.Shapes.Item(2).TextFrame.TextRange

for ($i = 1; $i -le $code.Paragraphs().Count; $i++) {
    $para = $code.Paragraphs($i)
    $t = $para.Text.TrimEnd([char]13, [char]10)
    if ($t -eq ";Author: Capt Jeff Falkinburg, USAF") {
        $start = $t.IndexOf("Jeff Falkinburg, ") + 1
        $len = "Jeff Falkinburg, ".Length
        $placeholder = "X".PadRight($len, "X")
        $sub = $para.Characters($start, $len)
        $sub.Text = $placeholder
        $sub2 = $para.Characters($start, $len)
        $sub2.Text = "Jeff Falkinburg, "
    }
}

# ---------------------------------------------------------------------------
# 3. Insert the new ABI slide right before the last slide ("Lab 2
#    Introduction"), so it becomes the new 16th slide.
# ---------------------------------------------------------------------------
$insertAt = $p.Slides.Count
$newSlide = $p.Slides.Add($insertAt, 16)

$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Application Binary Interface (ABI"
$title.Font.Bold = $true
$closeParen = $title.InsertAfter(")")
$closeParen.Font.Bold = $true

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Convention "
$r = $body.InsertAfter("of specifying which registers are used for arguments passed in to a subroutine and which are used to pass back results. ")
$r = $body.InsertAfter([char]13 + "For subroutines in the MSP430 use ")
$r = $body.InsertAfter("r12, r13, r14, and r15 to pass arguments to your subroutine. ")
$r = $body.InsertAfter([char]13 + "Use the stack if you have more than ")
$r = $body.InsertAfter("four arguments")
$r = $body.InsertAfter([char]13)
